$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing F/G cell values (daily antigen test counts revised)
$ws.Range("F250").Value = 17656
$ws.Range("F267").Value = 15446
$ws.Range("G267").Value = 867
$ws.Range("F268").Value = 17529
$ws.Range("F271").Value = 45844
$ws.Range("G271").Value = 1707
$ws.Range("F272").Value = 30120
$ws.Range("G272").Value = 1527
$ws.Range("F273").Value = 31827
$ws.Range("G273").Value = 1632
$ws.Range("F274").Value = 28609
$ws.Range("F275").Value = 30469
$ws.Range("F276").Value = 11491
$ws.Range("F278").Value = 31057
$ws.Range("F279").Value = 42461
$ws.Range("F281").Value = 46799
$ws.Range("F282").Value = 46408
$ws.Range("F283").Value = 17267
$ws.Range("F284").Value = 1209
$ws.Range("F285").Value = 42532
$ws.Range("F286").Value = 55322
$ws.Range("F288").Value = 59312
$ws.Range("F289").Value = 63068
$ws.Range("F290").Value = 17669
$ws.Range("F291").Value = 15200
$ws.Range("F292").Value = 83130
$ws.Range("F293").Value = 83345
$ws.Range("F295").Value = 17323
$ws.Range("F296").Value = 1882
$ws.Range("F297").Value = 2437
$ws.Range("F298").Value = 3279
$ws.Range("F299").Value = 66286
$ws.Range("F300").Value = 73072
$ws.Range("F301").Value = 72393
$ws.Range("F302").Value = 78857
$ws.Range("F303").Value = 9716
$ws.Range("F304").Value = 6156
$ws.Range("F305").Value = 3402
$ws.Range("F306").Value = 76467
$ws.Range("F307").Value = 75606
$ws.Range("F308").Value = 15726
$ws.Range("F309").Value = 78337
$ws.Range("G309").Value = 5566
$ws.Range("F310").Value = 79395
$ws.Range("F311").Value = 61552
$ws.Range("F312").Value = 28246
$ws.Range("F313").Value = 76610
$ws.Range("F314").Value = 65436
$ws.Range("F315").Value = 56564
$ws.Range("F316").Value = 50892
$ws.Range("F317").Value = 63828
$ws.Range("F318").Value = 49445
$ws.Range("F319").Value = 41413
$ws.Range("F320").Value = 73770
$ws.Range("F321").Value = 90908
$ws.Range("F322").Value = 110075
$ws.Range("F323").Value = 217637
$ws.Range("F324").Value = 250166
$ws.Range("F325").Value = 775362
$ws.Range("F326").Value = 418274
$ws.Range("F327").Value = 224306
$ws.Range("F328").Value = 181233
$ws.Range("F329").Value = 73451
$ws.Range("F330").Value = 71592
$ws.Range("F331").Value = 154322
$ws.Range("G331").Value = 2711
$ws.Range("F332").Value = 485587
$ws.Range("F333").Value = 255025
$ws.Range("F334").Value = 193319
$ws.Range("F335").Value = 150549
$ws.Range("F336").Value = 81974
$ws.Range("G336").Value = 2572
$ws.Range("F337").Value = 105086
$ws.Range("F338").Value = 221409
$ws.Range("F339").Value = 662675
$ws.Range("F340").Value = 387220
$ws.Range("F341").Value = 283499
$ws.Range("F342").Value = 178894
$ws.Range("F343").Value = 134101
$ws.Range("G343").Value = 2985
$ws.Range("F344").Value = 135987
$ws.Range("G344").Value = 2485
$ws.Range("F345").Value = 292311
$ws.Range("F346").Value = 675018
$ws.Range("F347").Value = 346643
$ws.Range("F348").Value = 232900
$ws.Range("F349").Value = 159844
$ws.Range("F350").Value = 127548
$ws.Range("F351").Value = 150522
$ws.Range("F352").Value = 307535
$ws.Range("F353").Value = 723567
$ws.Range("F354").Value = 316773
$ws.Range("F355").Value = 222113
$ws.Range("F356").Value = 160488
$ws.Range("G356").Value = 2885
$ws.Range("F357").Value = 138354
$ws.Range("F358").Value = 159027
$ws.Range("F359").Value = 321192
$ws.Range("F360").Value = 749744
$ws.Range("F361").Value = 332885
$ws.Range("F362").Value = 228996
$ws.Range("G362").Value = 3187
$ws.Range("F363").Value = 188780
$ws.Range("F364").Value = 168740
$ws.Range("G364").Value = 2490
$ws.Range("F365").Value = 184699
$ws.Range("G365").Value = 2399
$ws.Range("F366").Value = 339518
$ws.Range("G366").Value = 2844
$ws.Range("F367").Value = 767142
$ws.Range("F368").Value = 346279
$ws.Range("F369").Value = 235084
$ws.Range("G369").Value = 2606
$ws.Range("F370").Value = 180992
$ws.Range("G370").Value = 2046
$ws.Range("F371").Value = 160247
$ws.Range("G371").Value = 1967
$ws.Range("F372").Value = 178521
$ws.Range("G372").Value = 1856
$ws.Range("F373").Value = 350215
$ws.Range("G373").Value = 2382
$ws.Range("F374").Value = 773767
$ws.Range("G374").Value = 3425
$ws.Range("F375").Value = 351321
$ws.Range("F376").Value = 222158
$ws.Range("F377").Value = 176634
$ws.Range("F378").Value = 157341
$ws.Range("F379").Value = 179752
$ws.Range("F380").Value = 344742
$ws.Range("F381").Value = 746789
$ws.Range("F383").Value = 221082
$ws.Range("F384").Value = 171901
$ws.Range("F385").Value = 150819
$ws.Range("F386").Value = 182826
$ws.Range("F387").Value = 351574
$ws.Range("F388").Value = 730971
$ws.Range("G388").Value = 2208
$ws.Range("F390").Value = 219928
$ws.Range("G390").Value = 1475
$ws.Range("F391").Value = 177586
$ws.Range("F392").Value = 221506
$ws.Range("F398").Value = 298834
$ws.Range("F399").Value = 201532
$ws.Range("F400").Value = 150272
$ws.Range("G400").Value = 768
$ws.Range("F401").Value = 272250
$ws.Range("F402").Value = 719671
$ws.Range("F403").Value = 351741
$ws.Range("F404").Value = 225040
$ws.Range("F405").Value = 173686
$ws.Range("F406").Value = 170484
$ws.Range("F407").Value = 157984
$ws.Range("F408").Value = 303946
$ws.Range("F409").Value = 706268
$ws.Range("F410").Value = 363481
$ws.Range("G410").Value = 633
$ws.Range("F411").Value = 225065
$ws.Range("F412").Value = 175963
$ws.Range("F413").Value = 149196
$ws.Range("F414").Value = 148667
$ws.Range("F415").Value = 306011
$ws.Range("F416").Value = 669752
$ws.Range("F422").Value = 297189
$ws.Range("F425").Value = 139122
$ws.Range("F427").Value = 90292
$ws.Range("F428").Value = 101952
$ws.Range("F432").Value = 123105
$ws.Range("F434").Value = 79639
$ws.Range("F435").Value = 81829
$ws.Range("G435").Value = 263
$ws.Range("F436").Value = 143201
$ws.Range("G436").Value = 348
$ws.Range("F437").Value = 164969
$ws.Range("F438").Value = 120514
$ws.Range("F439").Value = 87787
$ws.Range("G439").Value = 317
$ws.Range("F440").Value = 72322
$ws.Range("G440").Value = 212

# Add new row 441 for 20.05.2021 (date serial 44335)
$ws.Range("A441").Value = 44335
$ws.Range("B441").Value = 388391
$ws.Range("C441").Value = 5597
$ws.Range("D441").Value = 253
$ws.Range("E441").Value = 12272
$ws.Range("F441").Value = 52450
$ws.Range("G441").Value = 156

Write-Host "Applied covid daily stats update"
